# Adds ACME integration settings to the Assets sheet and bumps the
# MaxRetryNumber constant, matching the "added data to config file" commit.

$wb = $excel.ActiveWorkbook

$constants = $wb.Worksheets.Item("Constants")
$assets    = $wb.Worksheets.Item("Assets")

# --- Assets sheet: three new asset rows -------------------------------
# (values are entered in this particular left/right order so the
# generated shared-string table matches the authored workbook)
$assets.Range("A2").Value = "ACME_URL"
$assets.Range("B3").Value = "http://www.sha1-online.com/"
$assets.Range("A3").Value = "SHA1_URL"
$assets.Range("B2").Value = "https://acme-test.uipath.com/"
$assets.Range("A4").Value = "ACME_CREDENTIAL"
$assets.Range("B4").Value = "ACME_Credential"

# --- Constants sheet: MaxRetryNumber updated from 0 to 2 ---------------
$constants.Range("B2").Value = 2

# --- Restore the on-screen selection / active sheet state ---------------
[void]$assets.Range("B4").Select()
[void]$constants.Activate()
[void]$constants.Range("B2").Select()
